$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.055.40"
$ws.Range("E2").Value = "  +2.76%  "
$ws.Range("D3").Value = "2.952.66"
$ws.Range("E3").Value = "  +0.89%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.86"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.52"
$ws.Range("E6").Value = "  +2.55%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "2.950.10"
$ws.Range("E8").Value = "  +0.88%  "
$ws.Range("E9").Value = "  +1.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.13"
$ws.Range("E10").Value = "  +2.72%  "
$ws.Range("E11").Value = "  +6.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.441"
$ws.Range("E12").Value = "  +0.64%  "
$ws.Range("E13").Value = "  +4.81%  "
$ws.Range("E14").Value = "  -2.09%  "
$ws.Range("E15").Value = "  -0.61%  "
$ws.Range("D16").Value = "3.442.32"
$ws.Range("E16").Value = "  +0.90%  "
$ws.Range("D17").Value = "63.016.28"
$ws.Range("E17").Value = "  +2.81%  "
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("D19").Value = "2.946.36"
$ws.Range("E19").Value = "  +0.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "442.35"
$ws.Range("E20").Value = "  +2.42%  "
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("E22").Value = "  -0.99%  "
$ws.Range("E23").Value = "  -1.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.01"
$ws.Range("E24").Value = "  -1.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.11"
$ws.Range("E25").Value = "  +1.96%  "
$ws.Range("E26").Value = "  -2.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.78"
$ws.Range("E27").Value = "  +0.27%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.30"
$ws.Range("E29").Value = "  +5.95%  "
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("E31").Value = "  +0.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0000101"
$ws.Range("E32").Value = "  +15.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.47"
$ws.Range("E33").Value = "  -0.62%  "
$ws.Range("E34").Value = "  -1.01%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.992"
$ws.Range("E36").Value = "  -1.62%  "
$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.60"
$ws.Range("E37").Value = "  -0.43%  "
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.08"
$ws.Range("E38").Value = "  +3.71%  "
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.68"
$ws.Range("E39").Value = "  -0.35%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.04"
$ws.Range("E40").Value = "  +2.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.50"
$ws.Range("E41").Value = "  -0.67%  "
$ws.Range("E42").Value = "  -4.48%  "
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "38.75"
$ws.Range("E44").Value = "  -8.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "135.55"
$ws.Range("E45").Value = "  +1.42%  "
$ws.Range("D46").Value = "2.693.59"
$ws.Range("E46").Value = "  -0.15%  "
$ws.Range("E47").Value = "  -2.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "359.96"
$ws.Range("E48").Value = "  -1.65%  "
$ws.Range("E50").Value = "  -0.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.74"
